$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assign responsible person(s) to a few previously-unassigned open tasks.
$ws.Range("C4").Value = "Jesse, Jonas"
$ws.Range("C5").Value = "Jesse"
$ws.Range("C7").Value = "Jesse"
$ws.Range("C8").Value = "Jesse"

# "addOrder umschreiben" is now done -> switch status text and restyle the
# cell the same way the other "done" row (B2) is styled.
$ws.Range("B11").Value = "done"
$ws.Range("B2").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New TODO row for the logout-on-login-page regression.
$ws.Range("A12").Value = "Login zurückgehen in Login-Page verhindern"

# Selection follows the newly active cell, matching the saved view state.
[void]$ws.Range("A12").Select()
